# Turn the first row (currently "AREA" / "Ciudad") into a centered header
# row reading "Valor" / "Categoría", matching the uploaded workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:B1")

$ws.Range("A1").Value2 = "Valor"
$ws.Range("B1").Value2 = "Categoría"

$headerRange.HorizontalAlignment = -4108   # xlCenter
